# Browser changes and Linux driver added
# Update the VerifyDoctor sheet's Trail Start/End Date test data to new dates
# and move the active selection to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VerifyDoctor")

$ws.Range("C2").Value = "25/12/2021"
$ws.Range("D2").Value = "25/12/2021"

$ws.Range("C3").Value = "25/12/2022"
$ws.Range("D3").Value = "25/12/2021"

$ws.Range("C4").Value = "25/12/2021"
$ws.Range("D4").Value = "25/12/2021"

$ws.Range("C5").Value = "25/12/2021"

$ws.Range("D6").Value = "25/12/2021"

$ws.Activate()
$ws.Range("C4").Select()
